# Update Name of Algo - adjust imputed numeric results on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.430100000000002
$ws.Range("B9").Value = 6.688699999999999
$ws.Range("D12").Value = -7.314099999999999
$ws.Range("E15").Value = 16.2436
$ws.Range("B18").Value = 7.739199999999999
$ws.Range("B20").Value = 9.345999999999997
$ws.Range("D26").Value = -8.349200000000005
$ws.Range("B27").Value = 5.578600000000006
$ws.Range("D27").Value = -8.678099999999999
$ws.Range("D29").Value = -7.251500000000001
$ws.Range("D37").Value = -7.846500000000002
$ws.Range("D38").Value = -7.214699999999998
$ws.Range("E38").Value = 17.04129999999999
$ws.Range("E44").Value = 16.11789999999999
$ws.Range("D51").Value = -8.661100000000003
$ws.Range("E51").Value = 16.4661
$ws.Range("D55").Value = -8.624499999999998
$ws.Range("E57").Value = 16.5947
$ws.Range("E63").Value = 18.54510000000002
$ws.Range("B69").Value = 5.377399999999998
$ws.Range("D69").Value = -7.263599999999994
$ws.Range("D70").Value = -8.283999999999999
$ws.Range("E70").Value = 16.64259999999999
$ws.Range("B76").Value = 5.3881
$ws.Range("B82").Value = 6.622599999999995
$ws.Range("D83").Value = -8.949999999999996
$ws.Range("E99").Value = 16.58629999999999
$ws.Range("D102").Value = -7.810200000000002
